# LOM3043.xlsx: fix rows 10-26, which had question (col A) and
# answer (col B/C) text shifted out of alignment by one row; also
# fills in previously-missing answers and adds a trailing Requisitos row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Copy-CellFormat {
    param($srcAddr, $dstAddr)
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null
}

# Row 10
$ws.Range("B10").Value = 'Apresentar os princípios básicos da Seleção de Materiais para aplicação em Engenharia.'
$ws.Range("C10").Value = 'Apresentar os princípios básicos da Seleção de Materiais para aplicação em Engenharia.'
# Row 13
$ws.Range("A13").Clear()
$ws.Range("B13").Value = '5840622 - Miguel Justino Ribeiro Barboza'
$ws.Range("C13").Value = '5840622 - Miguel Justino Ribeiro Barboza'
# Row 14
$ws.Range("A14").Value = 'Programa resumido:'
Copy-CellFormat "B9" "B14"
$ws.Range("B14").Value = 'Aspectos gerais e critérios de seleção de materiais estruturais. Aspectos dos principais mecanismos de falha em componentes estruturais. Seleção de materiais e análise para diferentes modos de carregamento. Seleção de materiais sob diferentes condições de temperatura. Materiais resistentes à corrosão e oxidação. Tribologia: atrito e desgaste. Tratamentos superficiais.'
Copy-CellFormat "C9" "C14"
$ws.Range("C14").Value = 'Aspectos gerais e critérios de seleção de materiais estruturais. Aspectos dos principais mecanismos de falha em componentes estruturais. Seleção de materiais e análise para diferentes modos de carregamento. Seleção de materiais sob diferentes condições de temperatura. Materiais resistentes à corrosão e oxidação. Tribologia: atrito e desgaste. Tratamentos superficiais.'
# Row 15
$ws.Range("A15").Value = 'Short syllabus:'
$ws.Range("B15").Clear()
$ws.Range("C15").Clear()
# Row 16
$ws.Range("A16").Value = 'Programa:'
Copy-CellFormat "B9" "B16"
$ws.Range("B16").Value = '1. Principais mecanismos de falha em componentes estruturais: efeitos do meio e temperatura. Critérios de falha. 2. Seleção de materiais para aplicações sob a ação de cargas estáticas. Materiais metálicos, cerâmicos, poliméricos e compósitos. 3. Seleção de materiais para aplicações sob a ação de cargas dinâmicas: O fenômeno da fadiga e efeitos da presença de entalhes em componentes mecânicos. 4. Critérios de seleção de materiais para aplicações em temperaturas elevadas. O fenômeno da fluência e a tolerância ao dano. Seleção de materiais para alta temperatura. Aços especiais, superligas, materiais cerâmicos e compósitos. 5. Materiais para temperaturas criogênicas. A transição dúctil-frágil. 6. Aspectos fundamentais do estudo de tribologia: desgaste, atrito e tratamentos superficiais. 7. Fundamentos, seleção e proteção contra oxidação. 8. Seleção de materiais em meios corrosivos. Corrosão sob tensão.'
Copy-CellFormat "C9" "C16"
$ws.Range("C16").Value = '1. Principais mecanismos de falha em componentes estruturais: efeitos do meio e temperatura. Critérios de falha. 2. Seleção de materiais para aplicações sob a ação de cargas estáticas. Materiais metálicos, cerâmicos, poliméricos e compósitos. 3. Seleção de materiais para aplicações sob a ação de cargas dinâmicas: O fenômeno da fadiga e efeitos da presença de entalhes em componentes mecânicos. 4. Critérios de seleção de materiais para aplicações em temperaturas elevadas. O fenômeno da fluência e a tolerância ao dano. Seleção de materiais para alta temperatura. Aços especiais, superligas, materiais cerâmicos e compósitos. 5. Materiais para temperaturas criogênicas. A transição dúctil-frágil. 6. Aspectos fundamentais do estudo de tribologia: desgaste, atrito e tratamentos superficiais. 7. Fundamentos, seleção e proteção contra oxidação. 8. Seleção de materiais em meios corrosivos. Corrosão sob tensão.'
# Row 17
$ws.Range("A17").Value = 'Syllabus:'
# Row 18
$ws.Range("A18").Value = 'Avaliação:'
$ws.Range("B18").Clear()
$ws.Range("C18").Clear()
# Row 19
$ws.Range("A19").Value = 'Método:'
# Row 20
$ws.Range("A20").Value = 'Critério:'
# Row 21
$ws.Range("A21").Value = 'Norma de recuperação:'
# Row 22
$ws.Range("A22").Value = 'Bibliografia:'
Copy-CellFormat "B9" "B22"
$ws.Range("B22").Value = '1.Ashby, M. F. Materials Selection in Mechanical Design, Butterworth, Oxford, 2005. 2. ASM Metals Handbook - Properties and Selection: Irons, Steels and High - Performance Alloys - v.1 - 1990. 3. ASM Metals Handbook - Properties and Selection: Nonferrous Alloys and Special - Purpose Materials - v.2 - 1990. 4. Meyers, M.; Chawla, K. Mechanical Behavior of Materials. Ed. Cambridge University Press, 2009. 5. Van Vlack, L.H., Propriedades dos Materiais Cerâmicos. Ed. Edgard Blücher Ltda., 1973. 6. Dowling, E. M. Mechanical behavior of materials: engineering methods for deformation, fracture and fatigue. New Jersey, Prentice Hall, 1999. 7. Biasotto, E., Polímeros como Materiais de Engenharia. Ed. Edgard Blücher Ltda., 1991. 8. Rosen, S.L., Fundamental Principles of Polymeric Materials. Ed. John Wiley & Sons, Inc., 1993. 9. Bhushan, B. Introduction to Tribology, 2nd Edition, John Wiley & Sons. 2013. 10. Roberge, P. R. Corrosion engineering: principles and practice. The McGraw-Hill Companies, Inc., 2008. 11. Gentil, V. Corrosão, Ed. LTC, 2011. 12. Crane, F.A., Charles, J.A., Selection of Engineering Materials, Butterworth, 1984. 13. Chiaverini, V., Aços e Ferros Fundidos, Associação Brasileira de Materiais - ABM, São Paulo, 1988. 14. Reed, R. C. The superalloys: fundamentals and applications. Ed. Cambridge, USA, 2006.'
Copy-CellFormat "C9" "C22"
$ws.Range("C22").Value = '1.Ashby, M. F. Materials Selection in Mechanical Design, Butterworth, Oxford, 2005. 2. ASM Metals Handbook - Properties and Selection: Irons, Steels and High - Performance Alloys - v.1 - 1990. 3. ASM Metals Handbook - Properties and Selection: Nonferrous Alloys and Special - Purpose Materials - v.2 - 1990. 4. Meyers, M.; Chawla, K. Mechanical Behavior of Materials. Ed. Cambridge University Press, 2009. 5. Van Vlack, L.H., Propriedades dos Materiais Cerâmicos. Ed. Edgard Blücher Ltda., 1973. 6. Dowling, E. M. Mechanical behavior of materials: engineering methods for deformation, fracture and fatigue. New Jersey, Prentice Hall, 1999. 7. Biasotto, E., Polímeros como Materiais de Engenharia. Ed. Edgard Blücher Ltda., 1991. 8. Rosen, S.L., Fundamental Principles of Polymeric Materials. Ed. John Wiley & Sons, Inc., 1993. 9. Bhushan, B. Introduction to Tribology, 2nd Edition, John Wiley & Sons. 2013. 10. Roberge, P. R. Corrosion engineering: principles and practice. The McGraw-Hill Companies, Inc., 2008. 11. Gentil, V. Corrosão, Ed. LTC, 2011. 12. Crane, F.A., Charles, J.A., Selection of Engineering Materials, Butterworth, 1984. 13. Chiaverini, V., Aços e Ferros Fundidos, Associação Brasileira de Materiais - ABM, São Paulo, 1988. 14. Reed, R. C. The superalloys: fundamentals and applications. Ed. Cambridge, USA, 2006.'
# Row 23
Copy-CellFormat "A9" "A23"
$ws.Range("A23").Value = 'Requisitos:'
$ws.Range("B23").Clear()
$ws.Range("C23").Clear()
# Row 24
$ws.Range("B24").Value = 'LOM3036 -  Propriedades Mecânicas  (Requisito fraco)
'
$ws.Range("C24").Value = 'LOM3036 -  Propriedades Mecânicas  (Requisito fraco)
'
# Row 25
$ws.Range("B25").Value = 'LOM3057 -  Introdução aos Materiais Poliméricos  (Requisito fraco)
'
$ws.Range("C25").Value = 'LOM3057 -  Introdução aos Materiais Poliméricos  (Requisito fraco)
'
# Row 26
Copy-CellFormat "B9" "B26"
$ws.Range("B26").Value = 'LOM3082 -  Cerâmica Física  (Requisito fraco)
'
Copy-CellFormat "C9" "C26"
$ws.Range("C26").Value = 'LOM3082 -  Cerâmica Física  (Requisito fraco)
'

# --- Row heights: realign custom heights to match the corrected content ---
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).AutoFit()
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 120
$ws.Rows.Item(23).AutoFit()
$ws.Rows.Item(26).RowHeight = 30
